# 9th Stab - Cosmetic Changes
# Insert two new date columns (Jun_17, Jun_15) to the left of the existing
# data, shifting the previous Jun_13 / Jun_10 columns two places to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank columns before column B (old B/C -> D/E).
$ws.Range("B1:C1").EntireColumn.Insert()

# Give the two newly inserted columns the same width as the rest of the
# date columns (matches the original 8.0-character custom width).
$ws.Columns.Item(3).ColumnWidth = 7.14
$ws.Columns.Item(4).ColumnWidth = 7.14
$ws.Columns.Item(5).ColumnWidth = 7.14

# New header row values (most-recent-first ordering: Jun_17, Jun_15, Jun_13, Jun_10)
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Fill the new columns with the same placeholder rating text ("UN") used
# throughout the rest of the sheet for every data row.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}
